# "atualizei dados da bibi e add"
#
# A new daily-revenue record (dia 30, julho/2025) was added to the top of
# the "07/2025" block. Since the sheet is sorted newest-period-first with
# day-of-month ascending inside each period, the new record belongs right
# after the existing "dia 29" row for julho (row 30) and before the
# "06/2025" block (which started at row 31) — i.e. it becomes the new
# row 31, pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 31; everything from the old row 31 down
# (the 06/2025, 05/2025 and 04/2025 blocks) shifts down to make room.
$ws.Rows.Item(31).Insert()

# Fill in the new day's figures.
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 64446.8
$ws.Range("C31").Value = 7
$ws.Range("D31").Value = 2025
$ws.Range("E31").Value = "07/2025"
